$d = $word.ActiveDocument

# Locate the run text ")}" that closes the asTable(...) field call, e.g.
# "...'C3', 'F7')}". That text currently lives in a single run -- the
# parser (TokenIteratorFieldRewriterSplit) now expects the closing "}"
# of the field to live in its own run, separate from the ")" that closes
# the function call.
$rng = $d.Content
$found = $rng.Find.Execute(")}", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $matchStart = $rng.Start
    $matchEnd = $rng.End

    # Re-select the matched ")}" text as its own Range.
    $target = $d.Range($matchStart, $matchEnd)

    # Rewrite the matched text as two runs via OOXML so the trailing "}"
    # becomes a genuinely new run rather than just edited text inside the
    # existing one. The first run keeps the original run's rsid
    # (00120327) and formatting (en-US language); the second, newly
    # split-off run carries the same formatting but -- like any
    # brand-new run -- no rsid.
    $xml = '<?xml version="1.0" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" ' + `
        'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
        'xmlns:xml="http://www.w3.org/XML/1998/namespace">' + `
        '<w:body><w:p>' + `
        '<w:r w:rsidR="00120327"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' + `
        '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">}</w:t></w:r>' + `
        '</w:p></w:body></w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)
}
